$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (Fitness) values for rows 2 through 191 to 7293,
# matching the value already present in rows 192 onward.
$ws.Range("C2:C191").Value = 7293
